$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet "2022-Q3" right before "2022-Q2" (so it
#    becomes the second sheet overall, pushing every quarter sheet down by
#    one position).
# ---------------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($existingQ2)
$q3.Name = "2022-Q3"

# Header row (same headers used by every other quarterly fund-holdings sheet)
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Make the header row + the "A" index column bold / centered / bordered, to
# match the look of the other quarterly sheets.
$headerStyleRange = $q3.Range("B1:H1")
$headerStyleRange.Font.Bold = $true
$headerStyleRange.HorizontalAlignment = -4108
$headerStyleRange.VerticalAlignment = -4160
$headerStyleRange.Borders.LineStyle = 1

# The text-like columns (fund code, fund name, size, position, etc.) need to
# stay as plain text, exactly as authored, so force a text format before
# writing the values (otherwise leading zeros / trailing decimals would be
# reinterpreted as numbers).
$q3.Range("B2:G3").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("A2").Font.Bold = $true
$q3.Range("A2").HorizontalAlignment = -4108
$q3.Range("A2").VerticalAlignment = -4160
$q3.Range("A2").Borders.LineStyle = 1
$q3.Range("B2").Value = "009649"
$q3.Range("C2").Value = "嘉实精选平衡混合A"
$q3.Range("D2").Value = "0.07"
$q3.Range("E2").Value = "58.93"
$q3.Range("F2").Value = "2.23"
$q3.Range("G2").Value = "0.0016"
$q3.Range("H2").Value = 10

$q3.Range("A3").Value = 1
$q3.Range("A3").Font.Bold = $true
$q3.Range("A3").HorizontalAlignment = -4108
$q3.Range("A3").VerticalAlignment = -4160
$q3.Range("A3").Borders.LineStyle = 1
$q3.Range("B3").Value = "009650"
$q3.Range("C3").Value = "嘉实精选平衡混合C"
$q3.Range("D3").Value = "0.05"
$q3.Range("E3").Value = "58.93"
$q3.Range("F3").Value = "2.23"
$q3.Range("G3").Value = "0.0011"
$q3.Range("H3").Value = 10

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q3 at the
#    top of the data (row 2), pushing the existing rows down by one, and
#    renumber the "A" index column sequentially (0,1,2,...).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows("2:2").Insert()
$summary.Range("A2:D2").ClearFormats()

# Re-use the exact same cell style already used by the other index cells
# (column A) further down the sheet, instead of rebuilding it by hand.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0

# Renumber the index column for the rows that were pushed down (they kept
# their old index values after the insert, so fix them up sequentially).
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6

Write-Host "2022-Q3 sheet added and summary updated"
